$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.864.89"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.814.72"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4669"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3682"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07376"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8699"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.874.84"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07064"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.495"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008692"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "26.912.28"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.338"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "2.087.24"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "149.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.318"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08928"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7666"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.940"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.261"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5310"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.352"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1661"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.411"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4924"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.668"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06281"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
